$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.386.04"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.911.23"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.80%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.88"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4819"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4068"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08219"
$ws.Range("E9").Value = "  +2.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.021"
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.51"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.915.03"
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.038"
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.201"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.09"
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06801"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001037"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.71"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.007"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.402.60"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.632"
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.186"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.141.98"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.566"
$ws.Range("E26").Value = "  +10.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.28"
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.04"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.105"
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.31"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.019"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09553"
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.612"
$ws.Range("E33").Value = "  +4.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.549"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.365"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02283"
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06109"
$ws.Range("E37").Value = "  +1.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.180"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.056"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5970"
$ws.Range("E40").Value = "  +2.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.79"
$ws.Range("E41").Value = "  +7.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1847"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.279"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.396"
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07612"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.45"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5572"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.951"
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.81"
$ws.Range("E49").Value = "  +4.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.419"
$ws.Range("E50").Value = "  +3.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.20"
$ws.Range("E51").Value = "  +1.15%  "
